# WIP : exception for cumulative values. Added new line in excel for change in cash
#
# Insert a new row for the alternate "CashAndCashEquivalentsPeriodIncreaseDecrease"
# XBRL tag under the same "Change in Cash, Cash Equivalents" standardized label,
# right below the existing "Change in Cash" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 36 - formatting is carried down from row 35 above it,
# matching the existing "Change in Cash, Cash Equivalents" row.
$ws.Rows("36:36").Insert()

$ws.Range("A36").Value = "Change in Cash, Cash Equivalents"
$ws.Range("C36").Value = "Amount of increase (decrease) in cash and cash equivalents. Cash and cash equivalents are the amount of currency on hand as well as demand deposits with banks or financial institutions. Includes other kinds of accounts that have the general characteristics of demand deposits. Also includes short-term, highly liquid investments that are both readily convertible to known amounts of cash and so near their maturity that they present insignificant risk of changes in value because of changes in interest rates. Includes effect from exchange rate changes."
$ws.Range("B36").Value = "CashAndCashEquivalentsPeriodIncreaseDecrease"

$ws.Rows("36:36").RowHeight = 15

# Small formatting touch-up on what is now row 40 (previously row 39), matching
# the font used elsewhere in the blank rows below the table.
$ws.Range("B38").Copy()
$ws.Range("C40").PasteSpecial(-4122)

# Re-apply the AutoFilter over the expanded table range (was A1:C37, now A1:C38).
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:C38").AutoFilter(1)

# Keep the workbook's hidden _FilterDatabase defined name in sync with the filter.
$n = $wb.Names.Item(1)
$n.RefersTo = "=Sheet1!`$A`$1:`$C`$38"
